# CreateNewProduct.xlsx edit script
# Applies: column insert (Packsize/ProdQty before Color), header/data cell
# updates, clears stale G2/G3 values, sets B11, updates selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Packsize, ProdQty) before the existing Color
# column (O), pushing Color/PaintType/FinishType/blank from O:R to Q:T.
$ws.Columns("O:P").Insert()

# ---- Row 1 (headers) ----
$ws.Range("A1").Value = "PaintName"
$ws.Range("B1").Value = "ProductStd"
$ws.Range("C1").Value = "ProductCode"
$ws.Range("D1").Value = "HSNId"
$ws.Range("E1").Value = "Brand "
$ws.Range("F1").Value = "Category"
$ws.Range("G1").Value = "SubCategory"
$ws.Range("H1").Value = "SubChildCategory"
$ws.Range("I1").Value = "BasePaint1"
$ws.Range("J1").Value = "ProportionOfPaint1"
$ws.Range("K1").Value = "Tinter1"
$ws.Range("L1").Value = "ProportionOfTinter "
$ws.Range("M1").Value = "BasePaint2"
$ws.Range("N1").Value = "ProportionOfPaint2"
$ws.Range("O1").Value = "Packsize"
$ws.Range("P1").Value = "ProdQty"
$ws.Range("Q1").Value = "Color"
$ws.Range("R1").Value = "PaintType"
$ws.Range("S1").Value = "FinishType"

# ---- Row 2 (data) ----
$ws.Range("A2").Value = "Test New 110"
$ws.Range("B2").Value = "Standard Product"
$ws.Range("C2").Value = 320004
$ws.Range("D2").Value = 32099099
$ws.Range("E2").Value = "Jotun"
$ws.Range("F2").Value = "Tinters"
$ws.Range("G2").Value = ""
$ws.Range("I2").Value = "Test 123 (10 Litre)"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "Test 121 (10 litre)"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "Test 123 (10 litre)"
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = "Red"
$ws.Range("R2").Value = "Oil Based"
$ws.Range("S2").Value = "Matt"

# ---- Row 3 (data) ----
$ws.Range("A3").Value = "Test New 111"
$ws.Range("B3").Value = "Non Standard Product"
$ws.Range("C3").Value = 320005
$ws.Range("D3").Value = 32099099
$ws.Range("E3").Value = "Jotun"
$ws.Range("F3").Value = "Tinters"
$ws.Range("G3").Value = ""
$ws.Range("I3").Value = "Test 124 (10 Litre)"
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "Test 121 (10 litre)"
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = "Test 118 (10 litre)"
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 5
$ws.Range("Q3").Value = "Red"
$ws.Range("R3").Value = "Oil Based"
$ws.Range("S3").Value = "Matt"

# ---- New note cell ----
$ws.Range("B11").Value = "s"

# ---- Column width tweak on K (Tinter1), best achievable approximation ----
$ws.Columns("K:K").ColumnWidth = 15.43

# ---- Restore the sheet's extended used-range/dimension (A1:AM19) by
# lightly touching the far corner cells without leaving visible content ----
$ws.Cells.Item(1, 39).Font.Bold = $false
$ws.Cells.Item(19, 39).Font.Bold = $false

# ---- Selection / view state ----
$ws.Range("B11").Select()

Write-Output "edit complete"
